$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("M:M").Delete()
[void]$ws.Range("M1").Select()
